# Update "想去人数" (want-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1817
$ws1.Range("F4").Value  = 401
$ws1.Range("F5").Value  = 1486
$ws1.Range("F9").Value  = 13159
$ws1.Range("F10").Value = 13034
$ws1.Range("F11").Value = 990
$ws1.Range("F12").Value = 764
$ws1.Range("F17").Value = 2060
$ws1.Range("F20").Value = 37
$ws1.Range("F22").Value = 190

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 103

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1817
$ws4.Range("F5").Value  = 401
$ws4.Range("F6").Value  = 1486
$ws4.Range("F11").Value = 13159
$ws4.Range("F12").Value = 13034
$ws4.Range("F13").Value = 990
$ws4.Range("F14").Value = 764
$ws4.Range("F21").Value = 2060
$ws4.Range("F24").Value = 37
$ws4.Range("F28").Value = 190
$ws4.Range("F31").Value = 103
